# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Phoenix_Profits data across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1338.6177
$ws.Range("I98").Value = 1358.6333
$ws.Range("K98").Value = 1358.6333
$ws.Range("M98").Value = 139.3667

$ws.Range("H122").Value = 1338.6177
$ws.Range("I122").Value = 1358.6333
$ws.Range("K122").Value = 4075.8999
$ws.Range("M122").Value = -1625.8999

$ws.Range("H132").Value = 2793.35
$ws.Range("I132").Value = 2763.6924
$ws.Range("K132").Value = 8291.0772
$ws.Range("M132").Value = -5761.0772

$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -160140

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 3402.389
$ws.Range("I4").Value = 1481.5454
$ws.Range("J4").Value = 6420.857
$ws.Range("K4").Value = 1481.5454
$ws.Range("L4").Value = 6420.857
$ws.Range("M4").Value = -1365.5454
$ws.Range("N4").Value = -6652.857

$ws.Range("H61").Value = 3209.9744
$ws.Range("I61").Value = 2223.1538
$ws.Range("J61").Value = 5183.615
$ws.Range("K61").Value = 2223.1538
$ws.Range("L61").Value = 5183.615
$ws.Range("M61").Value = -2011.1538
$ws.Range("N61").Value = -5607.615

$ws.Range("H102").Value = 1511.9048
$ws.Range("I102").Value = 1458.9375
$ws.Range("K102").Value = 1458.9375
$ws.Range("M102").Value = 163.0625

$ws.Range("H110").Value = 2658.6943
$ws.Range("I110").Value = 2543.8965
$ws.Range("K110").Value = 2543.8965
$ws.Range("M110").Value = -498.8964999999998

$ws.Range("H136").Value = 3209.9744
$ws.Range("I136").Value = 2223.1538
$ws.Range("J136").Value = 5183.615
$ws.Range("K136").Value = 6669.4614
$ws.Range("L136").Value = 15550.845
$ws.Range("M136").Value = -4119.4614
$ws.Range("N136").Value = -20650.845

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1809.0952
$ws.Range("I20").Value = 1748.0667
$ws.Range("J20").Value = 1961.6666
$ws.Range("K20").Value = 1748.0667
$ws.Range("L20").Value = 1961.6666
$ws.Range("M20").Value = -1501.0667
$ws.Range("N20").Value = -2455.6666

$ws.Range("H94").Value = 9616118
$ws.Range("I94").Value = 11905326
$ws.Range("J94").Value = 1445
$ws.Range("K94").Value = 11905326
$ws.Range("L94").Value = 1445
$ws.Range("M94").Value = -11904875
$ws.Range("N94").Value = -2347

$ws.Range("H107").Value = 11465.679
$ws.Range("I107").Value = 13119.392
$ws.Range("K107").Value = 13119.392
$ws.Range("M107").Value = -11199.392

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1674900.1
$ws.Range("I4").Value = 5500.5
$ws.Range("J4").Value = 2509600
$ws.Range("K4").Value = 5500.5
$ws.Range("L4").Value = 2509600
$ws.Range("M4").Value = -5388.5
$ws.Range("N4").Value = -2509824

$ws.Range("H7").Value = 168.34483
$ws.Range("I7").Value = 122.44444
$ws.Range("J7").Value = 243.45454
$ws.Range("K7").Value = 122.44444
$ws.Range("L7").Value = 243.45454
$ws.Range("M7").Value = -9.44444
$ws.Range("N7").Value = -469.45454

$ws.Range("H31").Value = 2288.7407
$ws.Range("I31").Value = 1061.7894
$ws.Range("K31").Value = 1061.7894
$ws.Range("M31").Value = -766.7893999999999

$ws.Range("H34").Value = 2288.7407
$ws.Range("I34").Value = 1061.7894
$ws.Range("K34").Value = 1061.7894
$ws.Range("M34").Value = -859.7893999999999

$ws.Range("H39").Value = 18868.5
$ws.Range("I39").Value = 10589.6
$ws.Range("J39").Value = 32666.666
$ws.Range("K39").Value = 10589.6
$ws.Range("L39").Value = 32666.666
$ws.Range("M39").Value = -10198.6
$ws.Range("N39").Value = -33448.666

$ws.Range("H49").Value = 18868.5
$ws.Range("I49").Value = 10589.6
$ws.Range("J49").Value = 32666.666
$ws.Range("K49").Value = 10589.6
$ws.Range("L49").Value = 32666.666
$ws.Range("M49").Value = -10407.6
$ws.Range("N49").Value = -33030.666

$ws.Range("H58").Value = 3607.111
$ws.Range("J58").Value = 3092.5
$ws.Range("L58").Value = 3092.5
$ws.Range("N58").Value = -3498.5

$ws.Range("H99").Value = 4367.4707
$ws.Range("I99").Value = 3982
$ws.Range("J99").Value = 4801.125
$ws.Range("K99").Value = 3982
$ws.Range("L99").Value = 4801.125
$ws.Range("M99").Value = -2484
$ws.Range("N99").Value = -7797.125

$ws.Range("H122").Value = 7107.4707
$ws.Range("I122").Value = 7202.5835
$ws.Range("J122").Value = 6879.2
$ws.Range("K122").Value = 21607.7505
$ws.Range("L122").Value = 20637.6
$ws.Range("M122").Value = -19157.7505
$ws.Range("N122").Value = -25537.6

$ws.Range("H126").Value = 4367.4707
$ws.Range("I126").Value = 3982
$ws.Range("J126").Value = 4801.125
$ws.Range("K126").Value = 11946
$ws.Range("L126").Value = 14403.375
$ws.Range("M126").Value = -9476
$ws.Range("N126").Value = -19343.375

$ws.Range("H134").Value = 2202.1082
$ws.Range("I134").Value = 1846.9678
$ws.Range("J134").Value = 4037
$ws.Range("K134").Value = 5540.903399999999
$ws.Range("L134").Value = 12111
$ws.Range("M134").Value = -3005.903399999999
$ws.Range("N134").Value = -17181

$ws.Range("H136").Value = 3607.111
$ws.Range("J136").Value = 3092.5
$ws.Range("L136").Value = 9277.5
$ws.Range("N136").Value = -14377.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13427946
$ws.Range("I4").Value = 16767607
$ws.Range("K4").Value = 50302821
$ws.Range("M4").Value = -50302709

$ws.Range("H14").Value = 3299.95
$ws.Range("I14").Value = 3299.95
$ws.Range("K14").Value = 9899.849999999999
$ws.Range("M14").Value = -9726.849999999999

$ws.Range("H86").Value = 554.41174
$ws.Range("I86").Value = 355.30768
$ws.Range("K86").Value = 1065.92304
$ws.Range("M86").Value = 120.0769599999999

$ws.Range("H89").Value = 554.41174
$ws.Range("I89").Value = 355.30768
$ws.Range("K89").Value = 3197.76912
$ws.Range("M89").Value = 2730.23088

$ws.Range("H126").Value = 1273.25
$ws.Range("I126").Value = 1273.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3819.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 1120.25
$ws.Range("N126").ClearContents()

$ws.Range("H129").Value = 4272.0586
$ws.Range("J129").Value = 6018.091
$ws.Range("L129").Value = 18054.273
$ws.Range("N129").Value = -28054.273

$ws.Range("H130").Value = 7503.5
$ws.Range("J130").Value = 10032
$ws.Range("L130").Value = 30096
$ws.Range("N130").Value = -40136

$ws.Range("H136").Value = 3253.7856
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 83348970
$ws.Range("I14").Value = 142874620
$ws.Range("J14").Value = 13046.4
$ws.Range("K14").Value = 142874620
$ws.Range("L14").Value = 13046.4
$ws.Range("M14").Value = -142874452
$ws.Range("N14").Value = -13382.4

$ws.Range("H47").Value = 30497.6
$ws.Range("J47").Value = 30497.6
$ws.Range("L47").Value = 30497.6
$ws.Range("N47").Value = -31633.6

$ws.Range("H55").Value = 31128.445
$ws.Range("J55").Value = 31693.334
$ws.Range("L55").Value = 31693.334
$ws.Range("N55").Value = -32347.334

$ws.Range("H97").Value = 1215
$ws.Range("I97").Value = 1282.8
$ws.Range("K97").Value = 1282.8
$ws.Range("M97").Value = -786.8

$ws.Range("H132").Value = 10000
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3741.4285
$ws.Range("J22").Value = 4148.5
$ws.Range("L22").Value = 4148.5
$ws.Range("N22").Value = -4738.5

$ws.Range("H27").Value = 3741.4285
$ws.Range("J27").Value = 4148.5
$ws.Range("L27").Value = 4148.5
$ws.Range("N27").Value = -4362.5

$ws.Range("H40").Value = 5015.696
$ws.Range("I40").Value = 4159.1665
$ws.Range("K40").Value = 4159.1665
$ws.Range("M40").Value = -4023.1665

$ws.Range("H48").Value = 35022.5

$ws.Range("H55").Value = 1744.1904
$ws.Range("I55").Value = 388.5
$ws.Range("K55").Value = 388.5
$ws.Range("M55").Value = -215.5

$ws.Range("H82").Value = 2171.25
$ws.Range("I82").Value = 2171.25
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2171.25
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1810.25
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 2171.25
$ws.Range("I85").Value = 2171.25
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2171.25
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -923.25
$ws.Range("N85").ClearContents()

$ws.Range("H93").Value = 2181.625
$ws.Range("I93").Value = 2131.2
$ws.Range("K93").Value = 2131.2
$ws.Range("M93").Value = -883.1999999999998

$ws.Range("H122").Value = 25101
$ws.Range("I122").Value = 25101
$ws.Range("K122").Value = 75303
$ws.Range("M122").Value = -72853

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3509.625
$ws.Range("I17").Value = 3834.7144
$ws.Range("J17").Value = 1234
$ws.Range("K17").Value = 3834.7144
$ws.Range("L17").Value = 1234
$ws.Range("M17").Value = -3662.7144
$ws.Range("N17").Value = -1578

$ws.Range("H107").Value = 2462.6667
$ws.Range("I107").Value = 1262.5
$ws.Range("J107").Value = 6063.1665
$ws.Range("K107").Value = 3787.5
$ws.Range("L107").Value = 18189.4995
$ws.Range("M107").Value = -1867.5
$ws.Range("N107").Value = -22029.4995

$ws.Range("H113").Value = 1015
$ws.Range("I113").Value = 1067.3334
$ws.Range("K113").Value = 3202.0002
$ws.Range("M113").Value = -1032.0002

$ws.Range("H126").Value = 39380544
$ws.Range("J126").Value = 3047.25
$ws.Range("L126").Value = 9141.75
$ws.Range("N126").Value = -14081.75
